$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 72, shifting existing rows 72:186 down to 73:187.
$ws.Rows("72:72").Insert()

# Populate the newly inserted row 72 with the new record's data.
$ws.Range("A72").Value = 3
$ws.Range("B72").Value = "Femacal de La Calera"
$ws.Range("C72").Value = "Coquimbo"
$ws.Range("D72").Value = 44477
$ws.Range("E72").Value = 5
$ws.Range("F72").Value = 100112039
$ws.Range("G72").Value = "Ciboulette"
$ws.Range("H72").Value = "Sin especificar"
$ws.Range("I72").Value = "Primera"
$ws.Range("J72").Value = 160
$ws.Range("K72").Value = 1500
$ws.Range("L72").Value = 1500
$ws.Range("M72").Value = 1500
$ws.Range("N72").Value = "`$/docena de atados"
$ws.Range("O72").Value = "Provincia de Quillota"
$ws.Range("P72").Value = 500
$ws.Range("Q72").Value = 3
$ws.Range("R72").Value = "Hortaliza"
